$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "2.16.840.1.113883.2.9.2.120.4.4.b0f3ffcf25ce2aafc7dc901e2febc51f43837f4ca0fe3b6d1b02194e9047b6db.b6dc1bd5c6^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Range("E2").Value = "2.16.840.1.113883.2.9.2.110.4.4^UAT_GTW_ID1721652251026"
$ws.Range("F2").Value = "22-07-2024:14:44:12"

$ws.Range("D3").Value = "2.16.840.1.113883.2.9.2.120.4.4.b0f3ffcf25ce2aafc7dc901e2febc51f43837f4ca0fe3b6d1b02194e9047b6db.e402b1f1de^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Range("E3").Value = "2.16.840.1.113883.2.9.2.110.4.4^UAT_GTW_ID1721652240625"
$ws.Range("F3").Value = "22-07-2024:14:44:02"
